# Auto-generated edit script
# Applies numeric value updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# worksheets per the scraped market-data refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 6184.125
$ws.Range("I40").Value2 = 1466.6666
$ws.Range("J40").Value2 = 7272.769
$ws.Range("K40").Value2 = 1466.6666
$ws.Range("L40").Value2 = 7272.769
$ws.Range("M40").Value2 = -1291.6666
$ws.Range("N40").Value2 = -7622.769
$ws.Range("H87").Value2 = 19999.908
$ws.Range("J87").Value2 = 19999.908
$ws.Range("L87").Value2 = 19999.908
$ws.Range("N87").Value2 = -22495.908
$ws.Range("H90").Value2 = 19999.908
$ws.Range("J90").Value2 = 19999.908
$ws.Range("L90").Value2 = 59999.724
$ws.Range("N90").Value2 = -72479.724
$ws.Range("H100").Value2 = 5443.2383
$ws.Range("I100").Value2 = 2890.5
$ws.Range("J100").Value2 = 7763.909
$ws.Range("K100").Value2 = 2890.5
$ws.Range("L100").Value2 = 7763.909
$ws.Range("M100").Value2 = -2349.5
$ws.Range("N100").Value2 = -8845.909
$ws.Range("H101").Value2 = 1714.25
$ws.Range("I101").Value2 = 588.1667
$ws.Range("K101").Value2 = 1764.5001
$ws.Range("M101").Value2 = -142.5001
$ws.Range("H103").Value2 = 1160.8572
$ws.Range("I103").Value2 = 784
$ws.Range("J103").Value2 = 1663.3334
$ws.Range("K103").Value2 = 2352
$ws.Range("L103").Value2 = 4990.0002
$ws.Range("M103").Value2 = -1766
$ws.Range("N103").Value2 = -6162.0002
$ws.Range("H138").Value2 = 10419694
$ws.Range("I138").Value2 = 1764.5
$ws.Range("J138").Value2 = 11907970
$ws.Range("K138").Value2 = 5293.5
$ws.Range("L138").Value2 = 35723910
$ws.Range("M138").Value2 = -153.5
$ws.Range("N138").Value2 = -35734190

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value2 = 1714.0312
$ws.Range("I88").Value2 = 1921.3334
$ws.Range("J88").Value2 = 1531.1177
$ws.Range("K88").Value2 = 1921.3334
$ws.Range("L88").Value2 = 1531.1177
$ws.Range("M88").Value2 = -1515.3334
$ws.Range("N88").Value2 = -2343.1177
$ws.Range("H91").Value2 = 1714.0312
$ws.Range("I91").Value2 = 1921.3334
$ws.Range("J91").Value2 = 1531.1177
$ws.Range("K91").Value2 = 1921.3334
$ws.Range("L91").Value2 = 1531.1177
$ws.Range("M91").Value2 = -517.3334
$ws.Range("N91").Value2 = -4339.1177
$ws.Range("H113").Value2 = 74900
$ws.Range("J113").Value2 = 74900
$ws.Range("L113").Value2 = 74900
$ws.Range("N113").Value2 = -83578
$ws.Range("H132").Value2 = 4095.24
$ws.Range("I132").Value2 = 3998.875
$ws.Range("J132").Value2 = 4480.7
$ws.Range("K132").Value2 = 11996.625
$ws.Range("L132").Value2 = 13442.1
$ws.Range("M132").Value2 = -9466.625
$ws.Range("N132").Value2 = -18502.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 6125.125
$ws.Range("I86").Value2 = 4244.222
$ws.Range("J86").Value2 = 8543.429
$ws.Range("K86").Value2 = 4244.222
$ws.Range("L86").Value2 = 8543.429
$ws.Range("M86").Value2 = -3121.222
$ws.Range("N86").Value2 = -10789.429
$ws.Range("H89").Value2 = 6125.125
$ws.Range("I89").Value2 = 4244.222
$ws.Range("J89").Value2 = 8543.429
$ws.Range("K89").Value2 = 21221.11
$ws.Range("L89").Value2 = 42717.145
$ws.Range("M89").Value2 = -15605.11
$ws.Range("N89").Value2 = -53949.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 4472.294
$ws.Range("I31").Value2 = 3752.75
$ws.Range("K31").Value2 = 3752.75
$ws.Range("M31").Value2 = -3457.75
$ws.Range("H34").Value2 = 4472.294
$ws.Range("I34").Value2 = 3752.75
$ws.Range("K34").Value2 = 3752.75
$ws.Range("M34").Value2 = -3550.75
$ws.Range("H59").Value2 = 699999.5
$ws.Range("I59").Value2 = 699999.5
$ws.Range("K59").Value2 = 699999.5
$ws.Range("M59").Value2 = -698854.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value2 = 47623476
$ws.Range("J103").Value2 = 66670868
$ws.Range("L103").Value2 = 200012604
$ws.Range("N103").Value2 = -200014362
$ws.Range("H122").Value2 = 1810.1428
$ws.Range("J122").Value2 = 1860.9333
$ws.Range("L122").Value2 = 16748.3997
$ws.Range("N122").Value2 = -21648.3997
$ws.Range("H128").Value2 = 459976.6
$ws.Range("I128").Value2 = 459976.6
$ws.Range("K128").Value2 = 1379929.8
$ws.Range("M128").Value2 = -1374949.8
$ws.Range("H129").Value2 = 1549.8572
$ws.Range("I129").Value2 = 816.25
$ws.Range("J129").Value2 = 2528
$ws.Range("K129").Value2 = 2448.75
$ws.Range("L129").Value2 = 7584
$ws.Range("M129").Value2 = 2551.25
$ws.Range("N129").Value2 = -17584
$ws.Range("H133").Value2 = 1000
$ws.Range("I133").Value2 = 1000
$ws.Range("J133").Value2 = 0
$ws.Range("K133").Value2 = 3000
$ws.Range("L133").Value2 = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value2 = 2060
$ws.Range("H137").Value2 = 2675.625
$ws.Range("I137").Value2 = 1521.2
$ws.Range("J137").Value2 = 4599.6665
$ws.Range("K137").Value2 = 4563.6
$ws.Range("L137").Value2 = 13798.9995
$ws.Range("M137").Value2 = 536.3999999999996
$ws.Range("N137").Value2 = -23998.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value2 = 295.92307
$ws.Range("J107").Value2 = 235
$ws.Range("L107").Value2 = 235
$ws.Range("N107").Value2 = -4075
$ws.Range("H113").Value2 = 6771.8125
$ws.Range("J113").Value2 = 13599.8
$ws.Range("L113").Value2 = 13599.8
$ws.Range("N113").Value2 = -17939.8
$ws.Range("H122").Value2 = 2730.45
$ws.Range("I122").Value2 = 2153.5881
$ws.Range("J122").Value2 = 5999.3335
$ws.Range("K122").Value2 = 6460.7643
$ws.Range("L122").Value2 = 17998.0005
$ws.Range("M122").Value2 = -4010.7643
$ws.Range("N122").Value2 = -22898.0005
$ws.Range("H126").Value2 = 5516.7334
$ws.Range("I126").Value2 = 3676.7144
$ws.Range("J126").Value2 = 7126.75
$ws.Range("K126").Value2 = 11030.1432
$ws.Range("L126").Value2 = 21380.25
$ws.Range("M126").Value2 = -8560.143199999999
$ws.Range("N126").Value2 = -26320.25
$ws.Range("H132").Value2 = 5513.591
$ws.Range("I132").Value2 = 5153.2666
$ws.Range("J132").Value2 = 6285.7144
$ws.Range("K132").Value2 = 15459.7998
$ws.Range("L132").Value2 = 18857.1432
$ws.Range("M132").Value2 = -12929.7998
$ws.Range("N132").Value2 = -23917.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value2 = 7555
$ws.Range("J2").Value2 = 0
$ws.Range("L2").Value2 = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value2 = 1271.2858
$ws.Range("I22").Value2 = 1149.8334
$ws.Range("K22").Value2 = 1149.8334
$ws.Range("M22").Value2 = -854.8334
$ws.Range("H27").Value2 = 1271.2858
$ws.Range("I27").Value2 = 1149.8334
$ws.Range("K27").Value2 = 1149.8334
$ws.Range("M27").Value2 = -1042.8334
$ws.Range("H93").Value2 = 324733.84
$ws.Range("I93").Value2 = 2010.04
$ws.Range("K93").Value2 = 2010.04
$ws.Range("M93").Value2 = -762.04
$ws.Range("H136").Value2 = 3981.5918
$ws.Range("I136").Value2 = 2048.6956
$ws.Range("J136").Value2 = 5691.4614
$ws.Range("K136").Value2 = 6146.0868
$ws.Range("L136").Value2 = 17074.3842
$ws.Range("M136").Value2 = -3596.0868
$ws.Range("N136").Value2 = -22174.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value2 = 1500
$ws.Range("I31").Value2 = 1000
$ws.Range("J31").Value2 = 2000
$ws.Range("K31").Value2 = 1000
$ws.Range("L31").Value2 = 2000
$ws.Range("M31").Value2 = -652
$ws.Range("N31").Value2 = -2696
$ws.Range("H81").Value2 = 3332.6667
$ws.Range("I81").Value2 = 2752.7693
$ws.Range("K81").Value2 = 5505.5386
$ws.Range("M81").Value2 = -4444.5386
$ws.Range("H84").Value2 = 3332.6667
$ws.Range("I84").Value2 = 2752.7693
$ws.Range("K84").Value2 = 27527.693
$ws.Range("M84").Value2 = -22223.693
